# Appointment form update:
#  - add a new "complete" (Name of Person Completing form) field
#  - reword the "Any notes about this Appointment?" question
#  - make several fields required (type_appoint, lab_test, date_appoint,
#    date_appoint1, complete)
#  - remove the "image1" note/image row from the form

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")

# 1. Insert a new row 47 (copy of row 46's formatting) for the new
#    "complete" field, right above the existing "type_appoint" row.
$ws.Rows("46:46").Copy()
$ws.Rows("47:47").Insert()

$ws.Range("A47").Value = "string"
$ws.Range("B47").Value = "complete"
$ws.Range("C47").Value = "Name of Person Completing form`n"
$ws.Range("G47").ClearContents()

# Match the "required" cell's shading to the rest of the shaded block
# (same look as the appearance cell next to it).
$ws.Range("G46").Copy()
$ws.Range("H47").PasteSpecial(-4122)
$ws.Range("H47").Value = "yes"

# Extend the yes/no list validation (previously only on D46) down to the
# newly inserted row.
$ws.Range("D46:D47").Validation.Delete()
$validation = $ws.Range("D46:D47").Validation
$validation.Add(3, 1, 1, '"yes,no"')
$validation.ShowInput = $false
$validation.ShowError = $false

# 2. type_appoint (now row 48) becomes required.
$ws.Range("A48").Copy()
$ws.Range("H48").PasteSpecial(-4122)
$ws.Range("H48").Value = "yes"

# 3. welcome (now row 49): reword the label and drop the "h1 blue"
#    appearance styling.
$ws.Range("C49").Value = "Notes about this Appointment?"
$ws.Range("G49").ClearContents()

# 4. lab_test (now row 50) becomes required.
$ws.Range("A50").Copy()
$ws.Range("H50").PasteSpecial(-4122)
$ws.Range("H50").Value = "yes"

# 5. date_appoint (now row 51) becomes required.
$ws.Range("A51").Copy()
$ws.Range("H51").PasteSpecial(-4122)
$ws.Range("H51").Value = "yes"

# 6. Remove the "image1" row (now row 52) entirely.
$ws.Rows("52:52").Delete()

# 7. date_appoint1 (now row 52, after the image row removal) becomes
#    required.
$ws.Range("A52").Copy()
$ws.Range("H52").PasteSpecial(-4122)
$ws.Range("H52").Value = "yes"
